$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the number format currently used by the last row's date cell
# (A66) - this is the distinguishing "last row" format that should move
# to the newly appended row.
$lastRowFormat = $ws.Range("A66").NumberFormat

# The previous last row (66) becomes a regular data row now that a new
# row is appended, so it should use the normal datetime format (matching
# the row above it, row 65).
$ws.Range("A66").NumberFormat = $ws.Range("A65").NumberFormat

# Append the new day's data as row 67.
$ws.Range("A67").Value = 45654
$ws.Range("B67").Value = 159
$ws.Range("C67").Value = 150
$ws.Range("D67").Value = 156

# Give the new last row the number format that used to belong to the old
# last row.
$ws.Range("A67").NumberFormat = $lastRowFormat
